$wb = $excel.ActiveWorkbook

$allData = @{}
$allData["n1_d40"] = @(
    @(0.0212155, 0.00553272, 0.0267482, 14.1),
    @(0.0227718, 0.00562025, 0.0283921, 13.6),
    @(0.0206139, 0.00473788, 0.0253518, 14.2),
    @(0.0210196, 0.00486491, 0.0258845, 14),
    @(0.0208091, 0.00479512, 0.0256042, 13.9),
    @(0.021666, 0.00616065, 0.0278267, 14),
    @(0.0213643, 0.00538354, 0.0267478, 14),
    @(0.0211277, 0.00536333, 0.0264911, 13.9),
    @(0.0213639, 0.00556403, 0.0269279, 13.9),
    @(0.0225553, 0.00639248, 0.0289478, 13.6),
    @(0.02145071, 0.005441491, 0.02689221000000001, 13.92),
)

$allData["n1_d60"] = @(
    @(0.0409223, 0.0101907, 0.051113, 26.6),
    @(0.0405839, 0.0106785, 0.0512624, 26.2),
    @(0.0428405, 0.0105669, 0.0534074, 26.3),
    @(0.0503368, 0.0116593, 0.0619961, 26.8),
    @(0.0395168, 0.00926767, 0.0487845, 25.1),
    @(0.0382896, 0.00856849, 0.046858, 25.9),
    @(0.0440698, 0.011463, 0.0555327, 26.1),
    @(0.0434486, 0.0111441, 0.0545927, 26),
    @(0.0398996, 0.010547, 0.0504467, 25.4),
    @(0.0417122, 0.0109787, 0.0526908, 25.9),
    @(0.04216200999999999, 0.010506436, 0.05266843, 26.03),
)

$allData["n1_d80"] = @(
    @(0.0673102, 0.0153914, 0.0827016, 39.8),
    @(0.0660911, 0.0152354, 0.0813265, 40.1),
    @(0.0713197, 0.0157933, 0.087113, 40.2),
    @(0.0698058, 0.0157359, 0.0855417, 40),
    @(0.0706732, 0.0166065, 0.0872798, 39.7),
    @(0.0692395, 0.0159004, 0.0851399, 40.5),
    @(0.0738304, 0.0162444, 0.0900748, 40.6),
    @(0.0737759, 0.0158808, 0.0896567, 40.8),
    @(0.0699207, 0.0160721, 0.0859928, 39.8),
    @(0.0688109, 0.0145976, 0.0834085, 40.1),
    @(0.07007774000000001, 0.01574578, 0.08582353, 40.16000000000001),
)

$allData["n1_d100"] = @(
    @(0.0864198, 0.0203437, 0.106763, 46.9),
    @(0.0847165, 0.0195095, 0.104226, 47.9),
    @(0.11262, 0.0218868, 0.134507, 47),
    @(0.0783596, 0.0196872, 0.0980468, 48),
    @(0.0755546, 0.0193434, 0.094898, 48.4),
    @(0.0853757, 0.0213178, 0.106694, 46.7),
    @(0.0883624, 0.0196897, 0.108052, 46.6),
    @(0.0764924, 0.0188136, 0.0953059, 48),
    @(0.0800889, 0.0198771, 0.099966, 46.8),
    @(0.0932211, 0.0225669, 0.115788, 47.9),
    @(0.08612110000000002, 0.02030357, 0.10642467, 47.42),
)

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $rows = $allData[$ws.Name]
    if ($null -eq $rows) { continue }

    # Insert two new blank columns at G; shifts old peak_memory_mb/total_clauses/total_variables
    # (G,H,I) to (I,J,K).
    $ws.Columns.Item(7).Insert()
    $ws.Columns.Item(7).Insert()

    # Update header row
    $ws.Range("F1").Value = "encoding_time_sec"
    $ws.Range("G1").Value = "solving_time_sec"
    $ws.Range("H1").Value = "total_time_sec"

    # Update data rows 2..12 with new measured encoding/solving/total times and peak memory
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $i + 2
        $row = $rows[$i]
        $ws.Cells.Item($r, 6).Value = $row[0]
        $ws.Cells.Item($r, 7).Value = $row[1]
        $ws.Cells.Item($r, 8).Value = $row[2]
        $ws.Cells.Item($r, 9).Value = $row[3]
    }

    # Clear now-stale total_clauses / total_variables data (columns J and K),
    # keeping the (shifted) header labels in J1 / K1 intact.
    $ws.Range("J2:K12").ClearContents()
}

